$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.113.98"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "1.651.73"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("E4").Value = "  -0.27%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.46"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +0.19%  "

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5216"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +0.08%  "

$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -0.28%  "

$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2647"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  +0.99%  "

$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06340"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +0.73%  "

$ws.Range("E10").Value = "  -0.50%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07686"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -1.62%  "

$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.625"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +3.30%  "

$ws.Range("D13").Value = "1.672.89"
$ws.Range("E13").Value = "  +1.30%  "

$ws.Range("D14").Value = "1.879.65"
$ws.Range("E14").Value = "  +0.07%  "

$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5595"
$ws.Range("D15").Style = $style

$ws.Range("D16").Value = "0.0₅8160"
$ws.Range("E16").Value = "  +1.92%  "

$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.39"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  +0.80%  "

$ws.Range("D18").Value = "26.110.15"
$ws.Range("E18").Value = "  +0.10%  "

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.628"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  -0.09%  "

$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.48"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +3.99%  "

$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "191.15"
$ws.Range("D22").Style = $style

$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.935"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("E24").Value = "  -0.29%  "

$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.15"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -1.06%  "

$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1190"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -1.14%  "

$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.220"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +0.65%  "

$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.96"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  +0.34%  "

$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.512"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +2.49%  "

$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05481"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -3.77%  "

$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.271"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  +0.32%  "

$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.450"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -1.03%  "

$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.362"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -0.12%  "

$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.559"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -1.89%  "

$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9493"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -0.09%  "

$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.786"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -0.58%  "

$ws.Range("E37").Value = "  -0.41%  "

$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5639"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -0.37%  "

$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01578"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -0.71%  "

$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.859"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -1.71%  "

$ws.Range("E41").Value = "  -0.23%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.028.10"
$ws.Range("E42").Value = "  -2.60%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8309"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -1.09%  "

$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.35"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -2.18%  "

$ws.Range("D45").Value = "1.792.20"
$ws.Range("E45").Value = "  +0.13%  "

$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.67"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +0.55%  "

$ws.Range("D47").Value = "0.0₈109"
$ws.Range("E47").Value = "  +5.66%  "

$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9998"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -0.80%  "

$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4338"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -1.26%  "

$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.993"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +0.42%  "

$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05176"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -3.50%  "
